$d = $word.ActiveDocument

# Replace the four <id>...</id> occurrences whose content needs the "a" prefix removed.
# (A fifth occurrence, fig_p104v_1, and the already-correct p104v_2 stay as separate runs
# or matching text and are left untouched by Find/Replace when no match is needed.)

$d.Content.Find.Execute("<id>p104v_a1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p104v_1</id>", 2)

$d.Content.Find.Execute("<id>p104v_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p104v_2</id>", 2)

$d.Content.Find.Execute("<id>p104v_a3</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p104v_3</id>", 2)

$d.Content.Find.Execute("<id>p104v_a4</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p104v_4</id>", 2)
